# Weekly fruit/vegetable price update: a new weekly record is inserted
# at the top of the "Repollo" (Macroferia Regional de Talca) data block,
# pushing the existing rows 559-588 down to 560-589.
#
# The new row reuses the (now shifted) row 560 values as its starting
# point - i.e. it is effectively a duplicate of the prior top record -
# except for the date (column D), which becomes 45267.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 559; existing rows 559:588 shift down to 560:589.
$ws.Rows.Item(559).Insert()

# Populate the new row 559 by duplicating the row now sitting at 560
# (the record that used to occupy row 559), then fix up the date.
$ws.Range("A560:R560").Copy($ws.Range("A559:R559"))
$ws.Range("D559").Value = 45267
